$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 2725
$ws.Range("I8").Value = 2725
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 8175
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -8036
# Row 28
$ws.Range("H28").Value = 9239.888999999999
$ws.Range("I28").Value = 1160.9231
$ws.Range("J28").Value = 30245.2
$ws.Range("K28").Value = 1160.9231
$ws.Range("L28").Value = 30245.2
$ws.Range("M28").Value = -675.9231
$ws.Range("N28").Value = -31215.2
# Row 132
$ws.Range("H132").Value = 2301.7454
$ws.Range("I132").Value = 2122.6938
$ws.Range("J132").Value = 3764
$ws.Range("K132").Value = 6368.0814
$ws.Range("L132").Value = 11292
$ws.Range("M132").Value = -3838.0814
$ws.Range("N132").Value = -16352
# Row 137
$ws.Range("H137").Value = 1985953.1
$ws.Range("I137").Value = 2977473
$ws.Range("J137").Value = 2913.2856
$ws.Range("K137").Value = 8932419
$ws.Range("L137").Value = 8739.856800000001
$ws.Range("M137").Value = -8929869
$ws.Range("N137").Value = -13839.8568
# Row 138
$ws.Range("H138").Value = 4657.8306
$ws.Range("I138").Value = 5243.75
$ws.Range("J138").Value = 4508.234
$ws.Range("K138").Value = 15731.25
$ws.Range("L138").Value = 13524.702
$ws.Range("M138").Value = -10591.25
$ws.Range("N138").Value = -23804.702

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 13348857
$ws.Range("I32").Value = 18194250
$ws.Range("J32").Value = 24025
$ws.Range("K32").Value = 18194250
$ws.Range("L32").Value = 24025
$ws.Range("M32").Value = -18193963
$ws.Range("N32").Value = -24599
# Row 45
$ws.Range("H45").Value = 1955.3529
$ws.Range("I45").Value = 999.5714
$ws.Range("J45").Value = 2624.4
$ws.Range("K45").Value = 999.5714
$ws.Range("L45").Value = 2624.4
$ws.Range("M45").Value = -622.5714
$ws.Range("N45").Value = -3378.4
# Row 74
$ws.Range("H74").Value = 8773374
$ws.Range("I74").Value = 851.5454999999999
$ws.Range("J74").Value = 38464988
$ws.Range("K74").Value = 851.5454999999999
$ws.Range("L74").Value = 38464988
$ws.Range("M74").Value = 22.45450000000005
$ws.Range("N74").Value = -38466736
# Row 77
$ws.Range("H77").Value = 8773374
$ws.Range("I77").Value = 851.5454999999999
$ws.Range("J77").Value = 38464988
$ws.Range("K77").Value = 4257.7275
$ws.Range("L77").Value = 192324940
$ws.Range("M77").Value = 110.2725
$ws.Range("N77").Value = -192333676
# Row 80
$ws.Range("H80").Value = 19284.715
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 19284.715
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 19284.715
$ws.Range("N80").Value = -21280.715
# Row 83
$ws.Range("H83").Value = 19284.715
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 19284.715
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 57854.145
$ws.Range("N83").Value = -67838.145
# Row 132
$ws.Range("H132").Value = 2203152.5
$ws.Range("I132").Value = 5769.6665
$ws.Range("J132").Value = 6997442.5
$ws.Range("K132").Value = 17308.9995
$ws.Range("L132").Value = 20992327.5
$ws.Range("M132").Value = -14778.9995
$ws.Range("N132").Value = -20997387.5

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 90911910
$ws.Range("I20").Value = 2975.6
$ws.Range("J20").Value = 166669360
$ws.Range("K20").Value = 2975.6
$ws.Range("L20").Value = 166669360
$ws.Range("M20").Value = -2728.6
$ws.Range("N20").Value = -166669854

$ws = $wb.Worksheets.Item("CRP")
# Row 28
$ws.Range("H28").Value = 51885.75
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 51885.75
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 51885.75
$ws.Range("N28").Value = -52375.75
# Row 31
$ws.Range("H31").Value = 7471.6577
$ws.Range("I31").Value = 3607.5715
$ws.Range("J31").Value = 8947.036
$ws.Range("K31").Value = 3607.5715
$ws.Range("L31").Value = 8947.036
$ws.Range("M31").Value = -3312.5715
$ws.Range("N31").Value = -9537.036
# Row 34
$ws.Range("H34").Value = 7471.6577
$ws.Range("I34").Value = 3607.5715
$ws.Range("J34").Value = 8947.036
$ws.Range("K34").Value = 3607.5715
$ws.Range("L34").Value = 8947.036
$ws.Range("M34").Value = -3405.5715
$ws.Range("N34").Value = -9351.036
# Row 122
$ws.Range("H122").Value = 2192.7097
$ws.Range("I122").Value = 1964
$ws.Range("J122").Value = 2407.125
$ws.Range("K122").Value = 5892
$ws.Range("L122").Value = 7221.375
$ws.Range("M122").Value = -3442
$ws.Range("N122").Value = -12121.375

$ws = $wb.Worksheets.Item("CUL")
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
# Row 132
$ws.Range("H132").Value = 2611.7556
$ws.Range("I132").Value = 2381.9312
$ws.Range("J132").Value = 3028.3125
$ws.Range("K132").Value = 21437.3808
$ws.Range("L132").Value = 27254.8125
$ws.Range("M132").Value = -18907.3808
$ws.Range("N132").Value = -32314.8125

$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 10086.956
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 10086.956
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 10086.956
$ws.Range("N15").Value = -10662.956
# Row 81
$ws.Range("H81").Value = 10086.956
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 10086.956
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 10086.956
$ws.Range("N81").Value = -12082.956
# Row 84
$ws.Range("H84").Value = 10086.956
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 10086.956
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 30260.868
$ws.Range("N84").Value = -40244.868
# Row 132
$ws.Range("H132").Value = 50008064
$ws.Range("I132").Value = 66675884
$ws.Range("J132").Value = 4602.4
$ws.Range("K132").Value = 200027652
$ws.Range("L132").Value = 13807.2
$ws.Range("M132").Value = -200025122
$ws.Range("N132").Value = -18867.2

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 6366.5
$ws.Range("I22").Value = 705
$ws.Range("J22").Value = 9969.272000000001
$ws.Range("K22").Value = 705
$ws.Range("L22").Value = 9969.272000000001
$ws.Range("M22").Value = -410
$ws.Range("N22").Value = -10559.272
# Row 27
$ws.Range("H27").Value = 6366.5
$ws.Range("I27").Value = 705
$ws.Range("J27").Value = 9969.272000000001
$ws.Range("K27").Value = 705
$ws.Range("L27").Value = 9969.272000000001
$ws.Range("M27").Value = -598
$ws.Range("N27").Value = -10183.272
# Row 68
$ws.Range("H68").Value = 1458.1464
$ws.Range("I68").Value = 1435.36
$ws.Range("J68").Value = 1493.75
$ws.Range("K68").Value = 1435.36
$ws.Range("L68").Value = 1493.75
$ws.Range("M68").Value = -686.3599999999999
$ws.Range("N68").Value = -2991.75
# Row 71
$ws.Range("H71").Value = 1458.1464
$ws.Range("I71").Value = 1435.36
$ws.Range("J71").Value = 1493.75
$ws.Range("K71").Value = 7176.799999999999
$ws.Range("L71").Value = 7468.75
$ws.Range("M71").Value = -3432.799999999999
$ws.Range("N71").Value = -14956.75
# Row 132
$ws.Range("H132").Value = 2840.5293
$ws.Range("I132").Value = 2090.1304
$ws.Range("J132").Value = 4409.5454
$ws.Range("K132").Value = 6270.3912
$ws.Range("L132").Value = 13228.6362
$ws.Range("M132").Value = -3740.3912
$ws.Range("N132").Value = -18288.6362
# Row 136
$ws.Range("H136").Value = 6945761.5
$ws.Range("I136").Value = 1254.9546
$ws.Range("J136").Value = 83335336
$ws.Range("K136").Value = 3764.8638
$ws.Range("L136").Value = 250006008
$ws.Range("M136").Value = -1214.8638
$ws.Range("N136").Value = -250011108

$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 1016.3333
$ws.Range("I100").Value = 958.6667
$ws.Range("J100").Value = 1054.7778
$ws.Range("K100").Value = 1917.3334
$ws.Range("L100").Value = 2109.5556
$ws.Range("M100").Value = -1376.3334
$ws.Range("N100").Value = -3191.5556
# Row 122
$ws.Range("H122").Value = 2993.9524
$ws.Range("I122").Value = 2545.6155
$ws.Range("J122").Value = 3722.5
$ws.Range("K122").Value = 7636.8465
$ws.Range("L122").Value = 11167.5
$ws.Range("M122").Value = -5186.8465
$ws.Range("N122").Value = -16067.5
# Row 124
$ws.Range("H124").Value = 39996.668
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 39996.668
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 39996.668
$ws.Range("N124").Value = -49816.668
# Row 132
$ws.Range("H132").Value = 9117608
$ws.Range("I132").Value = 3100.087
$ws.Range("J132").Value = 32410242
$ws.Range("K132").Value = 9300.261
$ws.Range("L132").Value = 97230726
$ws.Range("M132").Value = -6770.261
$ws.Range("N132").Value = -97235786
